# Weekly update: add two new "Fruta" (Mango) price rows for
# Agricola del Norte S.A. de Arica, inserted before the current row 229.
# This pushes the existing rows 229-248 down to 231-250.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 229 (old rows 229-248 become 231-250)
$ws.Rows("229:230").Insert()

# New row 229: Especial quality, 2023-04-05 (serial 45021)
$ws.Cells.Item(229, 1).Value = 1
$ws.Cells.Item(229, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(229, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(229, 4).Value = 45021
$ws.Cells.Item(229, 5).Value = 15
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100108
$ws.Cells.Item(229, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(229, 9).Value = 100108002
$ws.Cells.Item(229, 10).Value = "Mango"
$ws.Cells.Item(229, 11).Value = "Sin especificar"
$ws.Cells.Item(229, 12).Value = "Especial"
$ws.Cells.Item(229, 13).Value = 500
$ws.Cells.Item(229, 14).Value = 5000
$ws.Cells.Item(229, 15).Value = 5200
$ws.Cells.Item(229, 16).Value = 5100
$ws.Cells.Item(229, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(229, 18).Value = "Perú"
$ws.Cells.Item(229, 19).Value = 1275
$ws.Cells.Item(229, 20).Value = 4

# New row 230: Primera quality, same date 2023-04-05 (serial 45021)
$ws.Cells.Item(230, 1).Value = 1
$ws.Cells.Item(230, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(230, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(230, 4).Value = 45021
$ws.Cells.Item(230, 5).Value = 15
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100108
$ws.Cells.Item(230, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(230, 9).Value = 100108002
$ws.Cells.Item(230, 10).Value = "Mango"
$ws.Cells.Item(230, 11).Value = "Sin especificar"
$ws.Cells.Item(230, 12).Value = "Primera"
$ws.Cells.Item(230, 13).Value = 700
$ws.Cells.Item(230, 14).Value = 5000
$ws.Cells.Item(230, 15).Value = 5200
$ws.Cells.Item(230, 16).Value = 5129
$ws.Cells.Item(230, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(230, 18).Value = "Perú"
$ws.Cells.Item(230, 19).Value = 1282
$ws.Cells.Item(230, 20).Value = 4
